# Refresh the cryptocurrency price/volume table (GitHub Actions scheduled update).
# For D-column prices that are plain decimal numbers (e.g. "1.000", "240.80"),
# Excel's automatic type detection would otherwise coerce the text into a
# numeric value, so we briefly force a text number format, assign the value,
# then ClearFormats() to drop back to the cell's original (default/no style)
# formatting - matching how the source file stores these as plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.411.04'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.851.13'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6291'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07663'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2935'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07749'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').Value = '1.854.65'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.00001100'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.022'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.57'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '2.106.79'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.138'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').Value = '29.442.25'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.00'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.45'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.445'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.67'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1386'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.382'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.318'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.466'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05717'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.129'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.046'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.849'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.161'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7074'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.781'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').Value = '1.217.74'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.505'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9066'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.72'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.31'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000121'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.127'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4014'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.015'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.681'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1131'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.86%  '
